# Sync attendance_reports: fix "Recorded By" ordering in column G
# Swap the order of specific comma-separated entries (exact text replacements).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, dnasr281@gmail.com"             = "dnasr281@gmail.com, System"
    "System, system, backup@backdoor.com"    = "system, System, backup@backdoor.com"
    "admin@admin.com, System"                = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com"    = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
